$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Column H (rows 2-27) currently stores boolean FALSE values that are
# displayed as "FALSE" via a custom "TRUE";"TRUE";"FALSE" number format.
# Replace them with the literal text string "False" (Text-formatted cell,
# not a boolean) as described in the commit ("changing FALSE to False").
$rng = $ws.Range("H2:H27")
$rng.NumberFormat = "@"

# A direct Value/Formula assignment of the literal text "False" gets
# auto-sensed back into a Boolean by Excel's smart entry, same as typing
# it into the Name Box would. Enter it as a formula that evaluates to the
# text string, then flatten the formula to its cached value via
# copy / paste-special-values so the stored cell ends up as a genuine
# text (shared-string) literal instead of a boolean or a formula.
$rng.Formula = '="False"'
$rng.Copy()
$rng.PasteSpecial(-4163)   # xlPasteValues
$excel.CutCopyMode = $false

# Match the author's final selection/view state on the sheet.
$null = $ws.Range("H3:H27").Select()
